$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 108 - this shifts the existing rows
# 108..116 down to 109..117, preserving all of their data/formatting.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new weekly record.
$ws.Cells.Item(108, 1).Value2 = 6
$ws.Cells.Item(108, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(108, 3).Value2 = "Metropolitana"
$ws.Range("D108").Value2 = 44461
$ws.Cells.Item(108, 5).Value2 = 13
$ws.Cells.Item(108, 6).Value2 = 100112026
$ws.Cells.Item(108, 7).Value2 = "Haba"
$ws.Cells.Item(108, 8).Value2 = "Sin especificar"
$ws.Cells.Item(108, 9).Value2 = "Primera"
$ws.Cells.Item(108, 10).Value2 = 500
$ws.Cells.Item(108, 11).Value2 = 12000
$ws.Cells.Item(108, 12).Value2 = 14000
$ws.Cells.Item(108, 13).Value2 = 12920
$ws.Cells.Item(108, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(108, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(108, 16).Value2 = 517
$ws.Cells.Item(108, 17).Value2 = 25
$ws.Cells.Item(108, 18).Value2 = "Hortaliza"
